# Iraq League workbook update (30-03-2024 19:32)
# The underlying data rows for three match-pairs were re-ordered (swapped)
# while keeping the sequential "id" column (column A) untouched.
# Affected row pairs: (22,23), (73,74), (135,136)
# Swap columns B (2) through AC (29) between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

$rowPairs = @(
    @(22, 23),
    @(73, 74),
    @(135, 136)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
